$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 8, 9, 10, 11 in column H currently hold "JM" (leftover from earlier
# project stage). At time of delivery, these are updated to "-" like the
# other rows, making the "JM" shared string unused.
$ws.Range("H8").Value = "-"
$ws.Range("H9").Value = "-"
$ws.Range("H10").Value = "-"
$ws.Range("H11").Value = "-"

$ws.Range("H11").Select() | Out-Null
